$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (13.85546875 -> 14.85546875)
$ws.Columns.Item(1).ColumnWidth = 14

# New row 3
$ws.Cells.Item(3, 1).Value = 42600.829050925924
$ws.Cells.Item(3, 2).Value = "Bag"
$ws.Cells.Item(3, 3).Value = 3227
$ws.Cells.Item(3, 4).Value = 5623
$ws.Cells.Item(3, 5).Value = 620
$ws.Cells.Item(3, 6).Value = 103
$ws.Cells.Item(3, 7).Value = 56
$ws.Cells.Item(3, 8).Value = 63
$ws.Cells.Item(3, 9).Value = 34
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0

# New row 4
$ws.Cells.Item(4, 1).Value = 42600.881851851853
$ws.Cells.Item(4, 2).Value = "Bag"
$ws.Cells.Item(4, 3).Value = 9115
$ws.Cells.Item(4, 4).Value = 5630
$ws.Cells.Item(4, 5).Value = 621
$ws.Cells.Item(4, 6).Value = 103
$ws.Cells.Item(4, 7).Value = 56
$ws.Cells.Item(4, 8).Value = 63
$ws.Cells.Item(4, 9).Value = 34
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
